$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), formatted like the existing
# header row (copy H1's formatting, which carries style index 1: bold,
# bordered, centered/top-aligned).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data values for columns I (I0) and J (IF), rows 2-15.
$values = @{
    2  = @{ I = 1; J = 5 }
    3  = @{ I = 1; J = 4 }
    4  = @{ I = 1; J = 6 }
    5  = @{ I = 1; J = 6 }
    6  = @{ I = 1; J = 6 }
    7  = @{ I = 1; J = 6 }
    8  = @{ I = 1; J = 2 }
    9  = @{ I = 1; J = 5 }
    10 = @{ I = 1; J = 4 }
    11 = @{ I = 6; J = 6 }
    12 = @{ I = 8; J = 9 }
    13 = @{ I = 5; J = 6 }
    14 = @{ I = 7; J = 9 }
    15 = @{ I = 1; J = 3 }
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 9).Value = $values[$row].I
    $ws.Cells.Item($row, 10).Value = $values[$row].J
}
